$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last two rows of the old table (13 rows -> 11 rows) while
# keeping the formatting of the rows that remain.
$ws.Range("A12:C13").Delete()

# Header row
$ws.Range("A1").Value = "NBA, Friday 23rd Feb 2024"
$ws.Range("B1").Value = "Ballgorithm"
$ws.Range("C1").Value = "ESPN"

# Matchup / prediction rows
$ws.Range("A2").Value = "Toronto Raptors (20-36) vs Atlanta Hawks (24-31)"
$ws.Range("B2").Value = "Toronto Raptors (53.5%)"
$ws.Range("C2").Value = "Atlanta Hawks (67.4%)"

$ws.Range("A3").Value = "Cleveland Cavaliers (36-18) vs Philadelphia 76ers (32-23)"
$ws.Range("B3").Value = "Cleveland Cavaliers (65.52%)"
$ws.Range("C3").Value = "Cleveland Cavaliers (68.8%)"

$ws.Range("A4").Value = "Phoenix Suns (33-23) vs Houston Rockets (24-31)"
$ws.Range("B4").Value = "Houston Rockets (67.86%)"
$ws.Range("C4").Value = "Phoenix Suns (60.4%)"

$ws.Range("A5").Value = "Los Angeles Clippers (36-18) vs Memphis Grizzlies (20-36)"
$ws.Range("B5").Value = "Los Angeles Clippers (76.92%)"
$ws.Range("C5").Value = "Los Angeles Clippers (80.2%)"

$ws.Range("A6").Value = "Miami Heat (30-25) vs New Orleans Pelicans (34-22)"
$ws.Range("B6").Value = "New Orleans Pelicans (62.96%)"
$ws.Range("C6").Value = "New Orleans Pelicans (66.6%)"

$ws.Range("A7").Value = "Washington Wizards (9-46) vs Oklahoma City Thunder (38-17)"
$ws.Range("B7").Value = "Oklahoma City Thunder (78.57%)"
$ws.Range("C7").Value = "Oklahoma City Thunder (92.2%)"

$ws.Range("A8").Value = "Charlotte Hornets (14-41) vs Golden State Warriors (28-26)"
$ws.Range("B8").Value = "Golden State Warriors (51.72%)"
$ws.Range("C8").Value = "Golden State Warriors (91.7%)"

$ws.Range("A9").Value = "Milwaukee Bucks (35-21) vs Minnesota Timberwolves (39-16)"
$ws.Range("B9").Value = "Minnesota Timberwolves (79.17%)"
$ws.Range("C9").Value = "Minnesota Timberwolves (78.0%)"

$ws.Range("A10").Value = "Denver Nuggets (37-19) vs Portland Trail Blazers (15-39)"
$ws.Range("B10").Value = "Denver Nuggets (81.48%)"
$ws.Range("C10").Value = "Denver Nuggets (68.7%)"

$ws.Range("A11").Value = "San Antonio Spurs (11-45) vs Los Angeles Lakers (31-27)"
$ws.Range("B11").Value = "Los Angeles Lakers (68.97%)"
$ws.Range("C11").Value = "Los Angeles Lakers (79.1%)"

$ws.Range("C12").Select()
